# Add a new worksheet, "View Message Details Fields", at the end of the
# workbook and populate it with the plan for the "View Message Details"
# page (field names + notes), matching the layout/coloring used on the
# existing CreateMessageFields-style planning sheets.

$wb = $excel.ActiveWorkbook

# --- Add the new sheet after the last existing sheet ---------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "View Message Details Fields"

# Column B is a bit wider than the default so the second-column labels fit.
$ws.Columns.Item(2).ColumnWidth = 16.63

# --- Page title ------------------------------------------------------------
$ws.Range("A1").Value = "Fields to Go on View Message Details Page"

# --- Field / note values (styles applied afterwards, in bulk) -------------
$ws.Range("A6").Value  = "Message Create Date"
$ws.Range("A7").Value  = "Message Created Department"
$ws.Range("A8").Value  = "Message Created By"

$ws.Range("A10").Value = "Patient"
$ws.Range("C10").Value = "Any info about patient?"
$ws.Range("A11").Value = "Owner"
$ws.Range("A12").Value = "Owner Contact"

$ws.Range("A14").Value = "Type of Call"
$ws.Range("A15").Value = "Additional Questions"

$ws.Range("A17").Value = "Routing History: "
$ws.Range("A18").Value = "Time"
$ws.Range("B18").Value = "Routed From"
$ws.Range("C18").Value = "Routed To"

$ws.Range("A20").Value = "Message/Notes Table"
$ws.Range("A21").Value = "Who"
$ws.Range("B21").Value = "Who Department"
$ws.Range("C21").Value = "Time"
$ws.Range("D21").Value = "Message/Note"
$ws.Range("F21").Value = "Status"

$ws.Range("A26").Value = "Add Note Button"
$ws.Range("D26").Value = "You can only see this if it is claimed by you"
$ws.Range("A28").Value = "Route To Input Box"
$ws.Range("A29").Value = "Route Button"

$ws.Range("A31").Value = "Mark as Completed Button"

$ws.Range("A4").Value = "Claim Button for Unclaimed Messages/Unclaimed button for already claimed messages"

# --- Highlight styling -------------------------------------------------
# "Good" (green) first so it is registered before "Bad" (red).
$ws.Range("A6:F21").Style = "Good"

$ws.Range("A4:H4").Style = "Bad"
$ws.Range("A26:H30").Style = "Bad"

# --- Selection / active sheet --------------------------------------------
[void]$ws.Range("G20").Select()
[void]$ws.Activate()
